$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-04-06 Sunday" "2025-04-07 Monday"

Replace-Text "93÷4=23, 1" "76÷7=10, 6"
Replace-Text "98÷7=14, 0" "29÷9=3, 2"
Replace-Text "97÷9=10, 7" "33÷5=6, 3"
Replace-Text "40÷9=4, 4" "31÷8=3, 7"
Replace-Text "68÷2=34, 0" "31÷9=3, 4"

Replace-Text "55÷4=13, 3" "49÷4=12, 1"
Replace-Text "85÷9=9, 4" "34÷7=4, 6"
Replace-Text "38÷6=6, 2" "22÷3=7, 1"
Replace-Text "20÷2=10, 0" "84÷7=12, 0"
Replace-Text "81÷4=20, 1" "54÷4=13, 2"

Replace-Text "49÷3=16, 1" "33÷2=16, 1"
Replace-Text "31÷5=6, 1" "14÷4=3, 2"
Replace-Text "28÷2=14, 0" "86÷7=12, 2"
Replace-Text "41÷3=13, 2" "30÷2=15, 0"
Replace-Text "32÷4=8, 0" "69÷9=7, 6"

Replace-Text "69÷2=34, 1" "86÷4=21, 2"
Replace-Text "84÷3=28, 0" "75÷9=8, 3"
Replace-Text "51÷5=10, 1" "37÷6=6, 1"
Replace-Text "94÷2=47, 0" "96÷7=13, 5"
Replace-Text "93÷9=10, 3" "68÷4=17, 0"

Replace-Text "73÷2=36, 1" "92÷3=30, 2"
Replace-Text "67÷4=16, 3" "99÷3=33, 0"
Replace-Text "24÷6=4, 0" "78÷4=19, 2"
Replace-Text "28÷6=4, 4" "19÷6=3, 1"
Replace-Text "42÷3=14, 0" "56÷3=18, 2"

Write-Output "Done"
